$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at B, shifting the existing "Intent" (B1) and "Locale" (C1)
# headers right to C1/D1 (formatting + column width carry over automatically).
$ws.Range("B1").EntireColumn.Insert()

# The old "Question" header (A1) becomes "Answer".
$ws.Range("A1").Value = "Answer"

# The newly inserted column gets the "Category" header.
$ws.Range("B1").Value = "Category"

# Match the saved view state (active cell on row 2, ready for data entry).
$ws.Range("A2").Select()
